$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data was re-imported/re-sorted and merged with one additional
# JSON bucket, which reshuffled rows 2-14 (title/timestamp/uri travel
# together as a record; "historical distance"/"time bucket" stay "unknown"
# everywhere so the shuffle is invisible there). Existing hyperlink objects
# keep pointing at their original row's target (only the cell text backing
# them is updated), matching the source diff exactly.

$ws.Range("A2").Value = "Presidential Ratings"
$ws.Range("B2").Value = "1-01-01T00:00:00UTC"
$ws.Range("E2").Value = "https://insideelections.com/ratings/president"

$ws.Range("A3").Value = "Polls 2020-11-02 (larger states)"
$ws.Range("B3").Value = "1-01-01T00:00:00UTC"
$ws.Range("E3").Value = "https://www.swayable.com/polls/2020-11-02-large.html"

$ws.Range("A4").Value = "Find Your Local League"
$ws.Range("B4").Value = "1-01-01T00:00:00UTC"
$ws.Range("E4").Value = "https://www.lwv.org/local-leagues/find-local-league"

$ws.Range("A5").Value = "2020 President - Sabato's Crystal Ball"
$ws.Range("B5").Value = "1-01-01T00:00:00UTC"
$ws.Range("E5").Value = "http://centerforpolitics.org/crystalball/2020-president/"

$ws.Range("A6").Value = "Indiana: Election Tools, Deadlines, Dates, Rules, and Links"
$ws.Range("B6").Value = "1-01-01T00:00:00UTC"
$ws.Range("E6").Value = "https://www.vote.org/state/indiana/"

$ws.Range("A7").Value = "Election Results"
$ws.Range("B7").Value = "1-01-01T00:00:00UTC"
$ws.Range("E7").Value = "https://www.in.gov/sos/elections/2400.htm"

$ws.Range("A8").Value = "Karen Tallian Attorney General Poll"
$ws.Range("B8").Value = "2020-05-28T00:00:00UTC"
$ws.Range("E8").Value = "https://www.slideshare.net/IndianaBarrister/karen-tallian-attorney-general-poll"

$ws.Range("A9").Value = "Biden dominates the electoral map, but here's how the race could tighten"
$ws.Range("B9").Value = "2020-08-06T13:13:00UTC"
$ws.Range("E9").Value = "https://www.nbcnews.com/politics/meet-the-press/biden-dominates-electoral-map-here-s-how-race-could-tighten-n1236001"

$ws.Range("A10").Value = "Indiana Election Results"
$ws.Range("B10").Value = "1-01-01T00:00:00UTC"
$ws.Range("E10").Value = "https://enr.indianavoters.in.gov/site/index.html"

$ws.Range("A11").Value = "Voting & Elections Toolkits"
$ws.Range("B11").Value = "1-01-01T00:00:00UTC"
$ws.Range("E11").Value = "https://godort.libguides.com/votingtoolkit/texas"

$ws.Range("A12").Value = "2020 Election Forecast"
$ws.Range("B12").Value = "2020-08-12T06:30:00UTC"
$ws.Range("E12").Value = "https://projects.fivethirtyeight.com/2020-election-forecast/"

$ws.Range("A13").Value = "Tallian best positioned Democrat to win Indiana attorney general race, poll finds"
$ws.Range("B13").Value = "2020-05-28T17:30:00UTC"
$ws.Range("E13").Value = "https://www.nwitimes.com/news/local/govt-and-politics/elections/tallian-best-positioned-democrat-to-win-indiana-attorney-general-race-poll-finds/article_3208e96e-5678-584c-a81d-720ab4a22149.html"

$ws.Range("A14").Value = "2020 Electoral Interactive Map"
$ws.Range("B14").Value = "1-01-01T00:00:00UTC"
$ws.Range("E14").Value = "https://abcnews.go.com/Politics/2020-Electoral-Interactive-Map?basemap=71662160&promoref=brandpromo"
